$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H1").Value = "CategorizedProducts"
$ws.Range("A11").Value = "id"
$ws.Range("A11").Font.Bold = $true
$ws.Range("B11").Select()
